# Weekly update: a new price observation (row) is inserted into the
# "Poroto granado" sheet at row 59, pushing the existing rows 59-77 down
# to 60-78. The new row carries the latest ("current") observation, and
# the sheet's used range grows from A1:R77 to A1:R78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 59 - shifts rows 59..77 down to 60..78
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new observation
$ws.Range("A59").Value = 7
$ws.Range("B59").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C59").Value = "Ñuble"
$ws.Range("D59").Value = 44588
$ws.Range("E59").Value = 16
$ws.Range("F59").Value = 100112030
$ws.Range("G59").Value = "Poroto granado"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 120
$ws.Range("K59").Value = 23000
$ws.Range("L59").Value = 24000
$ws.Range("M59").Value = 23500
$ws.Range("N59").Value = "$/saco 25 kilos"
$ws.Range("O59").Value = "Provincia de Diguillín"
$ws.Range("P59").Value = 940
$ws.Range("Q59").Value = 25
$ws.Range("R59").Value = "Hortaliza"
